$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 25.28841266666667
$ws.Range("H2").Value = 75.865238
$ws.Range("I2").Value = 0.08258585054448338
$ws.Range("J2").Value = 0.08258585054448338
$ws.Range("M2").Value = 0.03927866666666666
$ws.Range("N2").Value = 0.117836
$ws.Range("O2").Value = 0.7432525340448212
$ws.Range("P2").Value = 0.7432525340448213
$ws.Range("Q2").Value = 0.9932951316631111
$ws.Range("R2").Value = 8.939656184968001
$ws.Range("S2").Value = 0.06138214269343414
$ws.Range("T2").Value = 0.06138214269343415

# Row 3
$ws.Range("G3").Value = 25.28841266666667
$ws.Range("H3").Value = 75.865238
$ws.Range("I3").Value = 0.08258585054448338
$ws.Range("J3").Value = 0.08258585054448338
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01356833333333333
$ws.Range("N3").Value = 0.040705
$ws.Range("O3").Value = 0.2567474659551788
$ws.Range("P3").Value = 0.2567474659551788
$ws.Range("Q3").Value = 0.3431216125322223
$ws.Range("R3").Value = 3.08809451279
$ws.Range("S3").Value = 0.02120370785104923
$ws.Range("T3").Value = 0.02120370785104923

# Row 4
$ws.Range("I4").Value = 0.6369880815661784
$ws.Range("J4").Value = 0.6369880815661784
$ws.Range("M4").Value = 0.03927866666666666
$ws.Range("N4").Value = 0.117836
$ws.Range("O4").Value = 0.7432525340448212
$ws.Range("P4").Value = 0.7432525340448213
$ws.Range("Q4").Value = 7.661326440009333
$ws.Range("R4").Value = 68.95193796008401
$ws.Range("S4").Value = 0.4734430057804113
$ws.Range("T4").Value = 0.4734430057804114

# Row 5
$ws.Range("I5").Value = 0.6369880815661784
$ws.Range("J5").Value = 0.6369880815661784
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01356833333333333
$ws.Range("N5").Value = 0.040705
$ws.Range("O5").Value = 0.2567474659551788
$ws.Range("P5").Value = 0.2567474659551788
$ws.Range("Q5").Value = 2.646511191321667
$ws.Range("R5").Value = 23.818600721895
$ws.Range("S5").Value = 0.163545075785767
$ws.Range("T5").Value = 0.163545075785767

# Row 6
$ws.Range("G6").Value = 55.14511
$ws.Range("H6").Value = 165.43533
$ws.Range("I6").Value = 0.180090615917626
$ws.Range("J6").Value = 0.180090615917626
$ws.Range("M6").Value = 0.03927866666666666
$ws.Range("N6").Value = 0.117836
$ws.Range("O6").Value = 0.7432525340448212
$ws.Range("P6").Value = 0.7432525340448213
$ws.Range("Q6").Value = 2.166026393986666
$ws.Range("R6").Value = 19.49423754588
$ws.Range("S6").Value = 0.1338528066384681
$ws.Range("T6").Value = 0.1338528066384682

# Row 7
$ws.Range("G7").Value = 55.14511
$ws.Range("H7").Value = 165.43533
$ws.Range("I7").Value = 0.180090615917626
$ws.Range("J7").Value = 0.180090615917626
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01356833333333333
$ws.Range("N7").Value = 0.040705
$ws.Range("O7").Value = 0.2567474659551788
$ws.Range("P7").Value = 0.2567474659551788
$ws.Range("Q7").Value = 0.7482272341833333
$ws.Range("R7").Value = 6.734045107649999
$ws.Range("S7").Value = 0.04623780927915786
$ws.Range("T7").Value = 0.04623780927915786

# Row 8
$ws.Range("G8").Value = 30.723475
$ws.Range("H8").Value = 92.170425
$ws.Range("I8").Value = 0.1003354519717122
$ws.Range("J8").Value = 0.1003354519717122
$ws.Range("M8").Value = 0.03927866666666666
$ws.Range("N8").Value = 0.117836
$ws.Range("O8").Value = 0.7432525340448212
$ws.Range("P8").Value = 0.7432525340448213
$ws.Range("Q8").Value = 1.206777133366666
$ws.Range("R8").Value = 10.8609942003
$ws.Range("S8").Value = 0.07457457893250752
$ws.Range("T8").Value = 0.07457457893250753

# Row 9
$ws.Range("G9").Value = 30.723475
$ws.Range("H9").Value = 92.170425
$ws.Range("I9").Value = 0.1003354519717122
$ws.Range("J9").Value = 0.1003354519717122
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.01356833333333333
$ws.Range("N9").Value = 0.040705
$ws.Range("O9").Value = 0.2567474659551788
$ws.Range("P9").Value = 0.2567474659551788
$ws.Range("Q9").Value = 0.4168663499583333
$ws.Range("R9").Value = 3.751797149625
$ws.Range("S9").Value = 0.02576087303920465
$ws.Range("T9").Value = 0.02576087303920465
